$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '27.142.55'
$ws.Cells.Item(2, 5).Value = '  -2.62%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.865.83'
$ws.Cells.Item(3, 5).Value = '  -2.19%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9999'
$ws.Cells.Item(4, 5).Value = '  -0.34%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '306.91'
$ws.Cells.Item(5, 5).Value = '  -2.03%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9995'
$ws.Cells.Item(6, 5).Value = '  -0.31%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5170'
$ws.Cells.Item(7, 5).Value = '  +3.33%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3752'
$ws.Cells.Item(8, 5).Value = '  -1.66%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07154'
$ws.Cells.Item(9, 5).Value = '  -1.74%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.8855'
$ws.Cells.Item(10, 5).Value = '  -2.51%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '20.68'
$ws.Cells.Item(11, 5).Value = '  -0.71%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.880.29'
$ws.Cells.Item(12, 5).Value = '  -1.54%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.07562'
$ws.Cells.Item(13, 5).Value = '  -1.32%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.334'
$ws.Cells.Item(14, 5).Value = '  -2.65%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '89.39'
$ws.Cells.Item(15, 5).Value = '  -2.49%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.0000'
$ws.Cells.Item(16, 5).Value = '  -0.40%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.000008550'
$ws.Cells.Item(17, 5).Value = '  -1.96%  '

$ws.Cells.Item(18, 5).Value = '  -2.56%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '1.000'

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '27.187.43'
$ws.Cells.Item(20, 5).Value = '  -2.55%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.029'
$ws.Cells.Item(21, 5).Value = '  -2.63%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.118.03'
$ws.Cells.Item(22, 5).Value = '  -1.75%  '

$ws.Cells.Item(23, 5).Value = '  -2.15%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '6.476'
$ws.Cells.Item(24, 5).Value = '  -1.96%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '150.94'
$ws.Cells.Item(25, 5).Value = '  -2.01%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.849'
$ws.Cells.Item(26, 5).Value = '  -1.85%  '

$ws.Cells.Item(27, 5).Value = '  -2.12%  '

$ws.Cells.Item(28, 5).Value = '  -4.01%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '112.80'
$ws.Cells.Item(29, 5).Value = '  -2.25%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.745'
$ws.Cells.Item(30, 5).Value = '  -3.32%  '

$ws.Cells.Item(31, 5).Value = '  +1.00%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.09017'
$ws.Cells.Item(32, 5).Value = '  +0.49%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.05163'
$ws.Cells.Item(33, 5).Value = '  -1.67%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '3.098'
$ws.Cells.Item(34, 5).Value = '  -3.66%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.7549'
$ws.Cells.Item(35, 5).Value = '  -1.43%  '

$ws.Cells.Item(37, 5).Value = '  -1.50%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.531'
$ws.Cells.Item(38, 5).Value = '  -0.75%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.024'
$ws.Cells.Item(39, 5).Value = '  +0.40%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.081'
$ws.Cells.Item(40, 5).Value = '  -1.38%  '

$ws.Cells.Item(41, 5).Value = '  -4.62%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '6.659'
$ws.Cells.Item(42, 5).Value = '  -4.24%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '115.37'
$ws.Cells.Item(43, 5).Value = '  +3.43%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '8.502'
$ws.Cells.Item(44, 5).Value = '  +0.17%  '

$ws.Cells.Item(45, 5).Value = '  -2.15%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.4665'
$ws.Cells.Item(46, 5).Value = '  -2.96%  '

$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.9995'
$ws.Cells.Item(47, 5).Value = '  -0.34%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '10.13'
$ws.Cells.Item(48, 5).Value = '  -4.42%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.573'
$ws.Cells.Item(49, 5).Value = '  -3.61%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '64.90'
$ws.Cells.Item(50, 5).Value = '  -3.78%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '36.36'
$ws.Cells.Item(51, 5).Value = '  -1.72%  '
